$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" column (R) is being added, mirroring the existing 2020
# column (Q) both for the year header (row 4) and the data point (row 5).
# Copy each source cell onto its new neighbor so the number format, font,
# borders and alignment all carry over, then overwrite with the new values.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 3.6

# Reflect the post-edit active selection.
$ws.Range("O9").Select()
